# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the prices in column D for rows 14-16
$ws.Range("D14").Value = 1266.597
$ws.Range("D15").Value = 1546.566
$ws.Range("D16").Value = 1817.002
